$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new working-hours entry was logged for 2014-03-21. Insert a fresh row at
# position 83 so the existing "sum" rows (old rows 83-85) shift down to
# 84-86, making room right after the last data row (row 82).
$ws.Rows.Item(83).Insert()

# The previous last entry (row 82, 2014-03-20) had its end time corrected
# from 19:15 to 21:00.
$ws.Range("E82").Value = 0.875

# Populate the newly inserted row 83 with the new log entry:
# 2014-03-21, 12:30 -> 13:15.
$ws.Range("A83").Value = 2014
$ws.Range("B83").Value = 3
$ws.Range("C83").Value = 21
$ws.Range("D83").Value = 0.52083333333333337
$ws.Range("E83").Value = 0.55208333333333337

# Give the new row the same "time spent" formulas used by the rows above it.
$ws.Range("F83").Formula = "=(E83-D83)*24*60"
$ws.Range("G83").Formula = "=F83/60"

# Match the author's final selection in the refreshed workbook.
$ws.Range("A84").Select()
